# Inserts a new weekly price record at row 234 (shifting the subsequent
# records down by one row) in the "Hortaliza, Vega Central Mapocho de
# Santiago - Berenjena" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 234..243 down to 235..244, creating a blank row 234.
$ws.Rows.Item(234).Insert()

# Fill in the new record.
$ws.Cells.Item(234, 1).Value = 9
$ws.Cells.Item(234, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(234, 3).Value = 'Metropolitana'
$ws.Cells.Item(234, 4).Value = 44610
$ws.Cells.Item(234, 5).Value = 13
$ws.Cells.Item(234, 6).Value = 100112001
$ws.Cells.Item(234, 7).Value = 'Berenjena'
$ws.Cells.Item(234, 8).Value = 'Sin especificar'
$ws.Cells.Item(234, 9).Value = 'Primera'
$ws.Cells.Item(234, 10).Value = 61
$ws.Cells.Item(234, 11).Value = 11000
$ws.Cells.Item(234, 12).Value = 13000
$ws.Cells.Item(234, 13).Value = 12016
$ws.Cells.Item(234, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(234, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(234, 16).Value = 200
$ws.Cells.Item(234, 17).Value = 60
$ws.Cells.Item(234, 18).Value = 'Hortaliza'
